# Minor changes to script and excel
# 1. Shorten the severity header labels on both sheets (drop the
#    " nedsættelse [<=]" / " nedsættelse [>=]" suffixes).
# 2. Swap which sheet/cell is the active selection: "Maend" was the
#    active tab (selection C13); now "Kvinder" is the active tab
#    (selection G2) and "Maend" keeps selection G1.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Kvinder")
$ws2 = $wb.Worksheets.Item("Maend")

# Update header row labels (write in E,D,C,F,G order so the shared-string
# table is rebuilt in the same order as the target workbook).
$ws1.Range("E1").Value = "Moderat"
$ws1.Range("D1").Value = "Svag"
$ws1.Range("C1").Value = "Ingen"
$ws1.Range("F1").Value = "Alvorlig"
$ws1.Range("G1").Value = "Dybtgående"

$ws2.Range("E1").Value = "Moderat"
$ws2.Range("D1").Value = "Svag"
$ws2.Range("C1").Value = "Ingen"
$ws2.Range("F1").Value = "Alvorlig"
$ws2.Range("G1").Value = "Dybtgående"

# Update selections / active sheet. Select on "Maend" first (it loses the
# active tab), finishing on "Kvinder" so it ends up the active tab.
$ws2.Activate()
[void]$ws2.Range("G1").Select()

$ws1.Activate()
[void]$ws1.Range("G2").Select()
